$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1300.1666
$ws.Range("I18").Value = 1360.2
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1360.2
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -1076.2
$ws.Range("N18").Value = -1568

$ws.Range("H33").Value = 235.46666
$ws.Range("I33").Value = 235.46666
$ws.Range("K33").Value = 235.46666
$ws.Range("M33").Value = -6.46665999999999

$ws.Range("H98").Value = 7409.864
$ws.Range("I98").Value = 5263.684
$ws.Range("J98").Value = 21002.334
$ws.Range("K98").Value = 5263.684
$ws.Range("L98").Value = 21002.334
$ws.Range("M98").Value = -3765.684
$ws.Range("N98").Value = -23998.334

$ws.Range("H112").Value = 1988.4546
$ws.Range("J112").Value = 2097.2666
$ws.Range("L112").Value = 6291.7998
$ws.Range("N112").Value = -8507.799800000001

$ws.Range("H122").Value = 7409.864
$ws.Range("I122").Value = 5263.684
$ws.Range("J122").Value = 21002.334
$ws.Range("K122").Value = 15791.052
$ws.Range("L122").Value = 63007.00199999999
$ws.Range("M122").Value = -13341.052
$ws.Range("N122").Value = -67907.00199999999

$ws.Range("H129").Value = 964.8868
$ws.Range("J129").Value = 1118.2927
$ws.Range("L129").Value = 3354.8781
$ws.Range("N129").Value = -13354.8781

$ws.Range("H138").Value = 2087106.6
$ws.Range("I138").Value = 5265217
$ws.Range("J138").Value = 4896.3965
$ws.Range("K138").Value = 15795651
$ws.Range("L138").Value = 14689.1895
$ws.Range("M138").Value = -15790511
$ws.Range("N138").Value = -24969.1895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1007.94116
$ws.Range("I2").Value = 848.3333
$ws.Range("J2").Value = 1187.5
$ws.Range("K2").Value = 848.3333
$ws.Range("L2").Value = 1187.5
$ws.Range("M2").Value = -735.3333
$ws.Range("N2").Value = -1413.5

$ws.Range("H9").Value = 40000

$ws.Range("H20").Value = 40000

$ws.Range("H45").Value = 846.5714
$ws.Range("I45").Value = 762.4
$ws.Range("J45").Value = 1057
$ws.Range("K45").Value = 762.4
$ws.Range("L45").Value = 1057
$ws.Range("M45").Value = -385.4
$ws.Range("N45").Value = -1811

$ws.Range("H74").Value = 1412.5714
$ws.Range("I74").Value = 1211.0667
$ws.Range("J74").Value = 1916.3334
$ws.Range("K74").Value = 1211.0667
$ws.Range("L74").Value = 1916.3334
$ws.Range("M74").Value = -337.0667000000001
$ws.Range("N74").Value = -3664.3334

$ws.Range("H77").Value = 1412.5714
$ws.Range("I77").Value = 1211.0667
$ws.Range("J77").Value = 1916.3334
$ws.Range("K77").Value = 6055.333500000001
$ws.Range("L77").Value = 9581.666999999999
$ws.Range("M77").Value = -1687.333500000001
$ws.Range("N77").Value = -18317.667

$ws.Range("H88").Value = 2796.1
$ws.Range("I88").Value = 2216.6667
$ws.Range("J88").Value = 3044.4285
$ws.Range("K88").Value = 2216.6667
$ws.Range("L88").Value = 3044.4285
$ws.Range("M88").Value = -1810.6667
$ws.Range("N88").Value = -3856.4285

$ws.Range("H91").Value = 2796.1
$ws.Range("I91").Value = 2216.6667
$ws.Range("J91").Value = 3044.4285
$ws.Range("K91").Value = 2216.6667
$ws.Range("L91").Value = 3044.4285
$ws.Range("M91").Value = -812.6667000000002
$ws.Range("N91").Value = -5852.4285

$ws.Range("H116").Value = 1007.94116
$ws.Range("I116").Value = 848.3333
$ws.Range("J116").Value = 1187.5
$ws.Range("K116").Value = 848.3333
$ws.Range("L116").Value = 1187.5
$ws.Range("M116").Value = 1445.6667
$ws.Range("N116").Value = -5775.5

$ws.Range("H122").Value = 3089.9333
$ws.Range("I122").Value = 3244.5833
$ws.Range("K122").Value = 9733.749899999999
$ws.Range("M122").Value = -7283.749899999999

$ws.Range("H132").Value = 436264.1
$ws.Range("I132").Value = 541613.7
$ws.Range("J132").Value = 3160.2222
$ws.Range("K132").Value = 1624841.1
$ws.Range("L132").Value = 9480.6666
$ws.Range("M132").Value = -1622311.1
$ws.Range("N132").Value = -14540.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1007.94116
$ws.Range("I3").Value = 848.3333
$ws.Range("J3").Value = 1187.5
$ws.Range("K3").Value = 848.3333
$ws.Range("L3").Value = 1187.5
$ws.Range("M3").Value = -734.3333
$ws.Range("N3").Value = -1415.5

$ws.Range("H86").Value = 41855.88
$ws.Range("I86").Value = 1750.5
$ws.Range("J86").Value = 68592.8
$ws.Range("K86").Value = 1750.5
$ws.Range("L86").Value = 68592.8
$ws.Range("M86").Value = -627.5
$ws.Range("N86").Value = -70838.8

$ws.Range("H89").Value = 41855.88
$ws.Range("I89").Value = 1750.5
$ws.Range("J89").Value = 68592.8
$ws.Range("K89").Value = 8752.5
$ws.Range("L89").Value = 342964
$ws.Range("M89").Value = -3136.5
$ws.Range("N89").Value = -354196

$ws.Range("H99").Value = 1650.125
$ws.Range("I99").Value = 1066.6666
$ws.Range("J99").Value = 2000.2
$ws.Range("K99").Value = 1066.6666
$ws.Range("L99").Value = 2000.2
$ws.Range("M99").Value = 431.3334
$ws.Range("N99").Value = -4996.2

$ws.Range("H107").Value = 25119.521
$ws.Range("I107").Value = 32576.412
$ws.Range("K107").Value = 32576.412
$ws.Range("M107").Value = -30656.412

$ws.Range("H134").Value = 314560.44
$ws.Range("I134").Value = 371988.12
$ws.Range("J134").Value = 4450.8
$ws.Range("K134").Value = 1115964.36
$ws.Range("L134").Value = 13352.4
$ws.Range("M134").Value = -1113429.36
$ws.Range("N134").Value = -18422.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 3000
$ws.Range("K25").Value = 3000
$ws.Range("M25").Value = -2826

$ws.Range("H86").Value = 3869.25
$ws.Range("I86").Value = 1350
$ws.Range("J86").Value = 6388.5
$ws.Range("K86").Value = 1350
$ws.Range("L86").Value = 6388.5
$ws.Range("M86").Value = -227
$ws.Range("N86").Value = -8634.5

$ws.Range("H89").Value = 3869.25
$ws.Range("I89").Value = 1350
$ws.Range("J89").Value = 6388.5
$ws.Range("K89").Value = 6750
$ws.Range("L89").Value = 31942.5
$ws.Range("M89").Value = -1134
$ws.Range("N89").Value = -43174.5

$ws.Range("H122").Value = 1344.8
$ws.Range("I122").Value = 1133.8182
$ws.Range("J122").Value = 1925
$ws.Range("K122").Value = 3401.4546
$ws.Range("L122").Value = 5775
$ws.Range("M122").Value = -951.4546
$ws.Range("N122").Value = -10675

$ws.Range("H123").Value = 65780
$ws.Range("J123").Value = 65780
$ws.Range("L123").Value = 65780
$ws.Range("N123").Value = -75580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 9491.666999999999
$ws.Range("I70").Value = 10590
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 31770
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -31455
$ws.Range("N70").Value = -12630

$ws.Range("H73").Value = 9491.666999999999
$ws.Range("I73").Value = 10590
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 31770
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -30678
$ws.Range("N73").Value = -14184

$ws.Range("H75").Value = 5022.2666
$ws.Range("J75").Value = 6308.091
$ws.Range("L75").Value = 18924.273
$ws.Range("N75").Value = -20920.273

$ws.Range("H78").Value = 5022.2666
$ws.Range("J78").Value = 6308.091
$ws.Range("L78").Value = 56772.819
$ws.Range("N78").Value = -66756.819

$ws.Range("H113").Value = 8321.357
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 8915.308000000001
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 26745.924
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -31085.924

$ws.Range("H125").Value = 3414.2856
$ws.Range("I125").Value = 2030
$ws.Range("J125").Value = 3968
$ws.Range("K125").Value = 6090
$ws.Range("L125").Value = 11904
$ws.Range("M125").Value = -1170
$ws.Range("N125").Value = -21744

$ws.Range("H129").Value = 1725895.4
$ws.Range("J129").Value = 1924967.9
$ws.Range("L129").Value = 5774903.699999999
$ws.Range("N129").Value = -5784903.699999999

$ws.Range("H131").Value = 20411938
$ws.Range("J131").Value = 28573440
$ws.Range("L131").Value = 85720320
$ws.Range("N131").Value = -85730400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70.75
$ws.Range("I2").Value = 49.11111
$ws.Range("J2").Value = 98.57143000000001
$ws.Range("K2").Value = 49.11111
$ws.Range("L2").Value = 98.57143000000001
$ws.Range("M2").Value = 63.88889
$ws.Range("N2").Value = -324.57143

$ws.Range("H80").Value = 3501.1
$ws.Range("I80").Value = 2917.5
$ws.Range("J80").Value = 4376.5
$ws.Range("K80").Value = 2917.5
$ws.Range("L80").Value = 4376.5
$ws.Range("M80").Value = -1919.5
$ws.Range("N80").Value = -6372.5

$ws.Range("H83").Value = 3501.1
$ws.Range("I83").Value = 2917.5
$ws.Range("J83").Value = 4376.5
$ws.Range("K83").Value = 14587.5
$ws.Range("L83").Value = 21882.5
$ws.Range("M83").Value = -9595.5
$ws.Range("N83").Value = -31866.5

$ws.Range("H140").Value = 46468.57
$ws.Range("J140").Value = 46468.57
$ws.Range("L140").Value = 46468.57
$ws.Range("N140").Value = -56828.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 11833.333
$ws.Range("I35").Value = 2750
$ws.Range("K35").Value = 2750
$ws.Range("M35").Value = -2414

$ws.Range("H100").Value = 8736.875
$ws.Range("I100").Value = 14473.75
$ws.Range("K100").Value = 14473.75
$ws.Range("M100").Value = -13932.75

$ws.Range("H132").Value = 7636.879
$ws.Range("I132").Value = 7604.52
$ws.Range("K132").Value = 22813.56
$ws.Range("M132").Value = -20283.56

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4222.1113
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4333.1665
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4333.1665
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5581.1665

$ws.Range("H65").Value = 4222.1113
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4333.1665
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 21665.8325
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -27905.8325

$ws.Range("H113").Value = 1293
$ws.Range("I113").Value = 1293
$ws.Range("K113").Value = 3879
$ws.Range("M113").Value = -1709

$ws.Range("H122").Value = 13889900
$ws.Range("I122").Value = 19231552
$ws.Range("J122").Value = 1603.6
$ws.Range("K122").Value = 57694656
$ws.Range("L122").Value = 4810.799999999999
$ws.Range("M122").Value = -57692206
$ws.Range("N122").Value = -9710.799999999999
